$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 27, shifting existing rows 27-31 down to 28-32.
$ws.Rows.Item(27).Insert()

# Copy the style (date format) used in column D down for the new row, same
# as the rest of the data rows.
$ws.Cells.Item(27, 4).Value = $ws.Cells.Item(28, 4).Value

# Populate the new row 27 with the data for the new record.
$ws.Cells.Item(27, 1).Value = 10
$ws.Cells.Item(27, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(27, 3).Value = "La Araucanía"
$ws.Cells.Item(27, 4).Value = (Get-Date -Year 2022 -Month 7 -Day 27 -Hour 0 -Minute 0 -Second 0).Date
$ws.Cells.Item(27, 5).Value = 9
$ws.Cells.Item(27, 6).Value = 100112042
$ws.Cells.Item(27, 7).Value = "Locoto"
$ws.Cells.Item(27, 8).Value = "Sin especificar"
$ws.Cells.Item(27, 9).Value = "Primera"
$ws.Cells.Item(27, 10).Value = 140
$ws.Cells.Item(27, 11).Value = 3300
$ws.Cells.Item(27, 12).Value = 3300
$ws.Cells.Item(27, 13).Value = 3300
$ws.Cells.Item(27, 14).Value = "$/kilo"
$ws.Cells.Item(27, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(27, 16).Value = 3300
$ws.Cells.Item(27, 17).Value = 1
$ws.Cells.Item(27, 18).Value = "Hortaliza"
